# Updated cryptos list - refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet.
#
# Column D values are stored as plain text in the workbook (e.g. "37.869.07",
# "227.25"). Several of the new values look like ordinary decimal numbers
# (e.g. "227.75"), and Excel's COM `.Value` setter auto-coerces such strings
# to numeric cells. To preserve the original text storage, column D is
# written via `.Formula` with a leading apostrophe (forces text) and the
# cell style is then reset to "Normal" so no stray number-format style is
# left behind. Column E values always contain "%"/spaces, so they are never
# mistaken for numbers and can be set directly via `.Value`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (or $null to leave unchanged), new Volume(1h)
$updates = @(
    @{ Row = 2;  D = "37.860.20"; E = "  -0.43%  " },
    @{ Row = 3;  D = "2.032.26";  E = "  -1.05%  " },
    @{ Row = 4;  D = $null;       E = "  -0.09%  " },
    @{ Row = 5;  D = "227.75";    E = "  -0.78%  " },
    @{ Row = 6;  D = "0.613";     E = "  -0.59%  " },
    @{ Row = 7;  D = "60.32";     E = "  +2.65%  " },
    @{ Row = 8;  D = $null;       E = "  -0.01%  " },
    @{ Row = 9;  D = $null;       E = "  -1.17%  " },
    @{ Row = 10; D = "0.0817";    E = "  +0.98%  " },
    @{ Row = 11; D = $null;       E = "  +0.61%  " },
    @{ Row = 12; D = "2.332.17";  E = "  -1.00%  " },
    @{ Row = 13; D = "14.54";     E = "  -0.90%  " },
    @{ Row = 14; D = "21.41";     E = "  +2.63%  " },
    @{ Row = 15; D = $null;       E = "  +1.00%  " },
    @{ Row = 16; D = $null;       E = "  -2.37%  " },
    @{ Row = 17; D = "2.048.91";  E = "  -0.23%  " },
    @{ Row = 18; D = "37.797.92"; E = "  -0.42%  " },
    @{ Row = 19; D = "69.84";     E = "  +0.19%  " },
    @{ Row = 20; D = $null;       E = "  -5.20%  " },
    @{ Row = 21; D = "0.0₃0827";  E = "  -0.99%  " },
    @{ Row = 22; D = "224.45";    E = "  -0.11%  " },
    @{ Row = 23; D = "0.998";     E = "  -0.16%  " },
    @{ Row = 24; D = "2.42";      E = "  -0.44%  " },
    @{ Row = 25; D = $null;       E = "  +0.52%  " },
    @{ Row = 26; D = "167.10";    E = "  +0.49%  " },
    @{ Row = 27; D = $null;       E = "  +0.06%  " },
    @{ Row = 28; D = $null;       E = "  -4.22%  " },
    @{ Row = 29; D = "18.90";     E = "  -0.67%  " },
    @{ Row = 30; D = $null;       E = "  -3.96%  " },
    @{ Row = 31; D = $null;       E = "  +0.82%  " },
    @{ Row = 32; D = $null;       E = "  +4.09%  " },
    @{ Row = 33; D = $null;       E = "  -2.78%  " },
    @{ Row = 34; D = $null;       E = "  -0.61%  " },
    @{ Row = 35; D = "4.50";      E = "  -2.01%  " },
    @{ Row = 36; D = "6.48";      E = "  +6.14%  " },
    @{ Row = 37; D = $null;       E = "  -2.55%  " },
    @{ Row = 38; D = $null;       E = "  -1.09%  " },
    @{ Row = 39; D = "0.999";     E = "  -0.24%  " },
    @{ Row = 40; D = "1.523.55";  E = "  +2.40%  " },
    @{ Row = 41; D = "17.13";     E = "  +3.42%  " },
    @{ Row = 42; D = $null;       E = "  +0.04%  " },
    @{ Row = 43; D = "96.07";     E = "  -1.06%  " },
    @{ Row = 44; D = $null;       E = "  -0.49%  " },
    @{ Row = 45; D = "0.0912";    E = "  -1.31%  " },
    @{ Row = 46; D = $null;       E = "  -1.53%  " },
    @{ Row = 47; D = $null;       E = "  -3.52%  " },
    @{ Row = 48; D = $null;       E = "  -0.89%  " },
    @{ Row = 49; D = $null;       E = "  -0.23%  " },
    @{ Row = 50; D = $null;       E = "  +0.65%  " },
    @{ Row = 51; D = "2.222.38";  E = "  -0.93%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Formula = "'" + $u.D
        $cell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
